{"js": "// Replace the bullet text \"- Maybe (S_2n)\" with \"- Maybe (S_nxS_n)\"\n// (keeps the existing bold/size/lang run formatting intact).\nconst body = context.document.body;\nconst results = body.search(\"- Maybe (S_2n)\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '- Maybe (S_2n)' not found in document body.\");\n}\n\nfor (const range of results.items) {\n  range.insertText(\"- Maybe (S_nxS_n)\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the bullet text \"- Maybe (S_2n)\" with \"- Maybe (S_nxS_n)\"\n# (formatting of the existing run - bold, size 28, en-US - is preserved).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"- Maybe (S_2n)\",   # FindText\n    $true,              # MatchCase\n    $false,             # MatchWholeWord\n    $false,             # MatchWildcards\n    $false,             # MatchSoundsLike\n    $false,             # MatchAllWordForms\n    $true,              # Forward\n    1,                  # Wrap            (wdFindContinue)\n    $false,             # Format\n    \"- Maybe (S_nxS_n)\",# ReplaceWith\n    2                   # Replace         (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw \"Target text '- Maybe (S_2n)' not found in document.\"\n}\n"}
